$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.710.44'
$ws.Range('E2').Value = '  +1.81%  '
$ws.Range('D3').Value = '1.637.90'
$ws.Range('E3').Value = '  +1.98%  '
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').Value = '212.92'
$ws.Range('E5').Value = '  +0.17%  '
$ws.Range('D6').Value = '0.493'
$ws.Range('E6').Value = '  +1.67%  '
$ws.Range('E7').Value = '  -0.21%  '
$ws.Range('D8').Value = '0.251'
$ws.Range('E8').Value = '  +1.68%  '
$ws.Range('E9').Value = '  +1.68%  '
$ws.Range('D10').Value = '18.99'
$ws.Range('E10').Value = '  +4.68%  '
$ws.Range('E11').Value = '  +2.66%  '
$ws.Range('D12').Value = '1.867.93'
$ws.Range('E12').Value = '  +2.11%  '
$ws.Range('D13').Value = '1.623.82'
$ws.Range('E13').Value = '  +1.06%  '
$ws.Range('D14').Value = '4.06'
$ws.Range('E14').Value = '  +1.16%  '
$ws.Range('E15').Value = '  +2.13%  '
$ws.Range('D16').Value = '26.701.18'
$ws.Range('E16').Value = '  +1.81%  '
$ws.Range('D17').Value = '63.02'
$ws.Range('E17').Value = '  +1.69%  '
$ws.Range('E18').Value = '  +1.90%  '
$ws.Range('D19').Value = '209.90'
$ws.Range('E19').Value = '  +4.57%  '
$ws.Range('E20').Value = '  -0.19%  '
$ws.Range('E21').Value = '  +0.98%  '
$ws.Range('D22').Value = '9.41'
$ws.Range('E22').Value = '  +1.35%  '
$ws.Range('E23').Value = '  +2.18%  '
$ws.Range('D24').Value = '1.94'
$ws.Range('E24').Value = '  +3.73%  '
$ws.Range('D25').Value = '146.49'
$ws.Range('E25').Value = '  +1.18%  '
$ws.Range('E26').Value = '  -0.22%  '
$ws.Range('E27').Value = '  -0.74%  '
$ws.Range('D28').Value = '6.73'
$ws.Range('E28').Value = '  +2.59%  '
$ws.Range('D29').Value = '15.40'
$ws.Range('E29').Value = '  +1.29%  '
$ws.Range('D30').Value = '0.0516'
$ws.Range('E30').Value = '  +5.24%  '
$ws.Range('E31').Value = '  -0.38%  '
$ws.Range('E32').Value = '  +1.00%  '
$ws.Range('E33').Value = '  +1.56%  '
$ws.Range('E34').Value = '  +0.99%  '
$ws.Range('D35').Value = '2.40'
$ws.Range('E35').Value = '  -0.18%  '
$ws.Range('D36').Value = '1.170.92'
$ws.Range('E36').Value = '  +0.98%  '
$ws.Range('E37').Value = '  -0.42%  '
$ws.Range('D38').Value = '0.809'
$ws.Range('E38').Value = '  +2.75%  '
$ws.Range('E39').Value = '  -0.18%  '
$ws.Range('E40').Value = '  +0.15%  '
$ws.Range('D41').Value = '0.503'
$ws.Range('E41').Value = '  +1.21%  '
$ws.Range('D42').Value = '0.794'
$ws.Range('E42').Value = '  +1.44%  '
$ws.Range('E43').Value = '  +1.75%  '
$ws.Range('D44').Value = '1.776.19'
$ws.Range('E44').Value = '  +2.07%  '
$ws.Range('D45').Value = '92.42'
$ws.Range('E45').Value = '  +0.80%  '
$ws.Range('E46').Value = '  +2.84%  '
$ws.Range('E47').Value = '  +10.27%  '
$ws.Range('E48').Value = '  +1.00%  '
$ws.Range('E49').Value = '  +1.35%  '
$ws.Range('E50').Value = '  +0.41%  '
$ws.Range('D51').Value = '7.53'
$ws.Range('E51').Value = '  +4.20%  '
